# Updates the worker account-status (Estado de Cuenta) table.
# Rewrites rows 16-29 (columns C-G) so that each of the 7 workers now has
# one row for period 1911 (rows 16-22) followed by one row for period 1912
# (rows 23-29), and updates a few "Valor Mora" / "Salario Basico" amounts
# for the newly-arranged rows (part 1 of the new statement data).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = "1017162179"
$ws.Range("D16").Value = "SANDRA JOHANA ACEVEDO VANEGAS"
$ws.Range("E16").Value = "1911"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116
$ws.Range("C17").Value = "71314293"
$ws.Range("D17").Value = "JOHN FABER MARTINEZ ALZATE"
$ws.Range("E17").Value = "1911"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 828116
$ws.Range("C18").Value = "43186853"
$ws.Range("D18").Value = "BIBIANA ECHEVERRI RAMIREZ"
$ws.Range("E18").Value = "1911"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 828116
$ws.Range("C19").Value = "1128433590"
$ws.Range("D19").Value = "NATALI GIRALDO VALENCIA"
$ws.Range("E19").Value = "1911"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 828116
$ws.Range("C20").Value = "43108510"
$ws.Range("D20").Value = "TATIANA LICED ARDILA AYA"
$ws.Range("E20").Value = "1911"
$ws.Range("F20").Value = 33125
$ws.Range("G20").Value = 828116
$ws.Range("C21").Value = "43622306"
$ws.Range("D21").Value = "SANDRA MILENA ACEVEDO MARIN"
$ws.Range("E21").Value = "1911"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 828116
$ws.Range("C22").Value = "43163272"
$ws.Range("D22").Value = "CLAUDIA ANDREA VAHOS RODRIGUEZ"
$ws.Range("E22").Value = "1911"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 100000
$ws.Range("C23").Value = "1017162179"
$ws.Range("D23").Value = "SANDRA JOHANA ACEVEDO VANEGAS"
$ws.Range("E23").Value = "1912"
$ws.Range("F23").Value = 33125
$ws.Range("G23").Value = 828116
$ws.Range("C24").Value = "71314293"
$ws.Range("D24").Value = "JOHN FABER MARTINEZ ALZATE"
$ws.Range("E24").Value = "1912"
$ws.Range("F24").Value = 33125
$ws.Range("G24").Value = 828116
$ws.Range("C25").Value = "43186853"
$ws.Range("D25").Value = "BIBIANA ECHEVERRI RAMIREZ"
$ws.Range("E25").Value = "1912"
$ws.Range("F25").Value = 33125
$ws.Range("G25").Value = 828116
$ws.Range("C26").Value = "1128433590"
$ws.Range("D26").Value = "NATALI GIRALDO VALENCIA"
$ws.Range("E26").Value = "1912"
$ws.Range("F26").Value = 33125
$ws.Range("G26").Value = 828116
$ws.Range("C27").Value = "43108510"
$ws.Range("D27").Value = "TATIANA LICED ARDILA AYA"
$ws.Range("E27").Value = "1912"
$ws.Range("F27").Value = 33125
$ws.Range("G27").Value = 828116
$ws.Range("C28").Value = "43622306"
$ws.Range("D28").Value = "SANDRA MILENA ACEVEDO MARIN"
$ws.Range("E28").Value = "1912"
$ws.Range("F28").Value = 40000
$ws.Range("G28").Value = 828116
$ws.Range("C29").Value = "43163272"
$ws.Range("D29").Value = "CLAUDIA ANDREA VAHOS RODRIGUEZ"
$ws.Range("E29").Value = "1912"
$ws.Range("F29").Value = 40000
$ws.Range("G29").Value = 100000
